$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 7
$ws.Range("G9").Value = 6
$ws.Range("G10").Value = 7
$ws.Range("G11").Value = 7
$ws.Range("G12").Value = 5
$ws.Range("G13").Value = 3
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 7
$ws.Range("G16").Value = 4
